$d = $word.ActiveDocument

# Remove pre-existing _GoBack bookmark first
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
Write-Host "Exists after delete: $($d.Bookmarks.Exists('_GoBack'))"

$rng = $d.Content
$found = $rng.Find.Execute("Profesional:")
$rng.Collapse(0)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = '<w:p ' + $ns + ' w:rsidR="00B3053B" w:rsidRDefault="00B3053B" w:rsidP="00B3053B"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Profesional:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>c</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>profesional</w:t></w:r><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p>'
$res = $rng.InsertXML($xml)
Write-Host "res=$res"
